# Commit: "modelling CO2 emissions and relevant emission tax"
#
# Adds a new node/commodity relationship row for CO2 emissions to the
# "rel_for_node_basic_structure" sheet (node = Emission_CO2_sink,
# commodity = emission_CO2), inserted as the first data row (row 3),
# pushing the existing data rows down by one. Also updates the saved
# UI selection state: the "rel_for_node_basic_structure" tab is no
# longer the active tab (its selection becomes the newly inserted
# row), and the "bus" tab becomes the active tab instead.

$wb = $excel.ActiveWorkbook

$wsRel = $wb.Worksheets.Item("rel_for_node_basic_structure")
$wsRel.Activate()

# Insert a brand-new row above the current row 3 (the first data row),
# shifting all the existing relationship rows down by one.
$wsRel.Rows.Item(3).Insert()

# Populate the newly inserted row with the same relationship-class
# columns used by every other row, plus the new node/commodity pair.
$wsRel.Cells.Item(3, 1).Value = "node__commodity"
$wsRel.Cells.Item(3, 2).Value = "node__stochastic_structure"
$wsRel.Cells.Item(3, 3).Value = "node__temporal_block"
$wsRel.Cells.Item(3, 4).Value = "Emission_CO2_sink"
$wsRel.Cells.Item(3, 5).Value = "emission_CO2"
$wsRel.Cells.Item(3, 6).Value = "default"
$wsRel.Cells.Item(3, 7).Value = "blk_t1"
$wsRel.Cells.Item(3, 8).Value = "blk_t2"

# Leave the new row selected on this sheet, matching the saved view
# state after the edit.
$wsRel.Range("A3:H3").Select()

# Switch the active tab to "bus", matching the workbook's saved view
# state after the edit.
$wsBus = $wb.Worksheets.Item("bus")
$wsBus.Activate()
